$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2 through 345 all hold the same date serial value (45189)
# which needs to be bumped forward by one day to 45190.
$ws.Range("C2:C345").Value = 45190
